$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.362.56'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '3.241.18'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.132'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.63'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.420'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("D12").Value = '3.819.16'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.137'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '68.365.78'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000169'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = '3.237.67'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.43'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '390.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.61'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.512'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000118'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.189'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.83%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '163.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.48'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.820'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.08'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.92%  '
$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0684'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '340.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.91%  '
$ws.Range("D46").Value = '2.590.86'
$ws.Range("E46").Value = '  -4.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0280'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.22%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '31.35'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("E51").Value = '  -1.55%  '
